$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text looks like a plain number (e.g. "21.12") ---
# Force these to remain TEXT cells (matching the source data feed, which
# stores prices/volumes as inline strings, not numeric values) by setting
# the cell to Text format before writing the value.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.64"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3680"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07233"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8607"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.12"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.625"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.379"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06887"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "80.66"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008851"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.17"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.185"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.00"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.78"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.26"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.217"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.884"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.88"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08928"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7422"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.158"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.419"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.799"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.116"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01921"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5081"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.763"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1642"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.431"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.274"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.80"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.39"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.651"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06286"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4556"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.805"

# --- Cells whose new text is unambiguous (already non-numeric text) ---
$ws.Range("D2").Value = "26.910.18"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.811.80"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +3.05%  "
$ws.Range("D11").Value = "2.043.40"
$ws.Range("E11").Value = "  +16.64%  "
$ws.Range("E12").Value = "  +4.40%  "
$ws.Range("E13").Value = "  +4.61%  "
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "26.953.07"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("E22").Value = "  +3.54%  "
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "2.291.67"
$ws.Range("E24").Value = "  +16.09%  "
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("E28").Value = "  +3.73%  "
$ws.Range("E29").Value = "  +16.26%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("E33").Value = "  +6.36%  "
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  +3.57%  "
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("E40").Value = "  +3.48%  "
$ws.Range("E41").Value = "  +9.67%  "
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("E43").Value = "  +5.52%  "
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E48").Value = "  +5.21%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("E51").Value = "  +5.49%  "
